# Apply updated "想去人数" (F) and "最低票价" (G) figures to the
# 展览 (Exhibitions) sheet and the 全部类型 (All types) sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value = 415
$wsExpo.Range("F7").Value = 544
$wsExpo.Range("F8").Value = 64
$wsExpo.Range("F9").Value = 6741
$wsExpo.Range("F10").Value = 153
$wsExpo.Range("F15").Value = 1083
$wsExpo.Range("F16").Value = 16099
$wsExpo.Range("F19").Value = 327
$wsExpo.Range("F20").Value = 176
$wsExpo.Range("F22").Value = 11308
$wsExpo.Range("F23").Value = 7
$wsExpo.Range("F24").Value = 916
$wsExpo.Range("F25").Value = 4443
$wsExpo.Range("F26").Value = 300
$wsExpo.Range("F28").Value = 42
$wsExpo.Range("F29").Value = 34
$wsExpo.Range("F31").Value = 138
$wsExpo.Range("F32").Value = 5215
$wsExpo.Range("G32").Value = "不可售"

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 415
$wsAll.Range("F7").Value = 544
$wsAll.Range("F9").Value = 64
$wsAll.Range("F10").Value = 6741
$wsAll.Range("F11").Value = 153
$wsAll.Range("F15").Value = 68
$wsAll.Range("F17").Value = 1083
$wsAll.Range("F18").Value = 16099
$wsAll.Range("F21").Value = 327
$wsAll.Range("F22").Value = 176
$wsAll.Range("F26").Value = 11308
$wsAll.Range("F27").Value = 7
$wsAll.Range("F28").Value = 916
$wsAll.Range("F29").Value = 4443
$wsAll.Range("F30").Value = 300
$wsAll.Range("F32").Value = 42
$wsAll.Range("F33").Value = 34
$wsAll.Range("F35").Value = 138
$wsAll.Range("F36").Value = 5215
$wsAll.Range("G36").Value = "不可售"
